# Merge Session into User: append sample login/session rows (15-20) to the
# "in" sheet, mirroring sign-up (userId, email, password, class) records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: userId -1 (stored as text, like the legacy rows), hell/hell@gmail.com
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "-1"
$ws.Range("C15").Value = "hell@gmail.com"
$ws.Range("D15").Value = "hell"
$ws.Range("E15").Value = "Student"
$ws.Range("F15").Value = $false
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

# --- Row 16: userId 0 (stored as text), helloworld1/helloworld1@gmail.com
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "0"
$ws.Range("C16").Value = "helloworld1@gmail.com"
$ws.Range("D16").Value = "helloworld1"
$ws.Range("E16").Value = "Student"
$ws.Range("F16").Value = $false
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""

# --- Row 17: userId -2 (real number), helloworld2/helloworld2@gmail.com
$ws.Range("B17").Value = -2
$ws.Range("C17").Value = "helloworld2@gmail.com"
$ws.Range("D17").Value = "helloworld2"
$ws.Range("E17").Value = "Student"
$ws.Range("F17").Value = $false
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""

# --- Row 18: userId -3 (real number), helloworld3/helloworld3@gmail.com
$ws.Range("B18").Value = -3
$ws.Range("C18").Value = "helloworld3@gmail.com"
$ws.Range("D18").Value = "helloworld3"
$ws.Range("E18").Value = "Student"
$ws.Range("F18").Value = $false
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""

# --- Row 19: userId -4 (real number), helloworld4/helloworld4@gmail.com
$ws.Range("B19").Value = -4
$ws.Range("C19").Value = "helloworld4@gmail.com"
$ws.Range("D19").Value = "helloworld4"
$ws.Range("E19").Value = "Student"
$ws.Range("F19").Value = $false
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = 0

# --- Row 20: userId -5 (real number), helloworld5/helloworld5@gmail.com
$ws.Range("B20").Value = -5
$ws.Range("C20").Value = "helloworld5@gmail.com"
$ws.Range("D20").Value = "helloworld5"
$ws.Range("E20").Value = "Student"
$ws.Range("F20").Value = $false
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = 0
